# daysofweek_list.xlsx edit:
#  - remove the "fra" language block (rows 9-15), shifting "ara" up to rows 9-15
#  - replace the (now shifted) "ara" block's lang_code/code/name with "hin" data
#  - leave day_seq / is_global_working / is_active columns untouched (values already match)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the fra rows (9-15), shifting ara rows (16-22) up into 9-15.
$ws.Range("A9:F15").EntireRow.Delete() | Out-Null

# 2) Overwrite the language code, numeric day code and day name for what is now
#    the hin block in rows 9-15.
$hinRows = @(
    @{ Row = 9;  Code = 101; Name = "रवि" },
    @{ Row = 10; Code = 102; Name = "सोमवार" },
    @{ Row = 11; Code = 103; Name = "मंगल" },
    @{ Row = 12; Code = 104; Name = "बुध" },
    @{ Row = 13; Code = 105; Name = "इकट्ठा करना" },
    @{ Row = 14; Code = 106; Name = "शुक्र" },
    @{ Row = 15; Code = 107; Name = "बैठा" }
)

foreach ($item in $hinRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = "hin"

    # Column B is formatted as Text ("@") at the column level, so a plain
    # Value assignment would store the number as a text string. Flip the
    # cell to General just long enough to write a true numeric value, then
    # restore the original (Text) number format so the cell's style/format
    # stays exactly as it was.
    $codeCell = $ws.Cells.Item($r, 2)
    $origFormat = $codeCell.NumberFormat
    $codeCell.NumberFormat = "General"
    $codeCell.Value = $item.Code
    $codeCell.NumberFormat = $origFormat

    $ws.Cells.Item($r, 3).Value = $item.Name
}
